# Updated Passwords to include FTP (ftp.womencoders.org)
# Adds a new "FTP Server" column (F) to the Sheet1 credentials table,
# reusing the "womencoders" username and adding the "peoplespace88" password.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell, styled like the other bold header cells (B1:E1)
$ws.Range("F1").Value = "FTP Server"
$ws.Range("F1").Font.Bold = $true

# New data rows for the FTP Server column
$ws.Range("F2").Value = "womencoders"
$ws.Range("F3").Value = "peoplespace88"

# Widen the new column to fit its contents
$ws.Columns.Item(6).ColumnWidth = 16.67

# Move the active selection to A4, below the table
$ws.Range("A4").Select()
